$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two source rows were dropped from the dataset ("RM 232" and "SC 92"),
# which shifts every subsequent row up. Delete the lower row first so the
# row number of the upper one doesn't change before it is removed.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the shift, a different set of "C" (column D) values is missing for
# the remaining SC rows - update them to match the new missing-data mask.
$ws.Range("D26").Value = -13.8
$ws.Range("D27").Value = ""
$ws.Range("D28").Value = ""
$ws.Range("D29").Value = -13
$ws.Range("D30").Value = -13.6
$ws.Range("D31").Value = ""
$ws.Range("D32").Value = ""
